# Update the five data rows (rows 2-6) of Sheet1 with refreshed sample
# values, per the "Update example Excel data files." commit.
#
# Columns A and G hold literal text in the source file (asset name, and a
# pre-formatted "NN.N%" string) rather than Excel-native types, so we force
# Text number-formatting before writing them and then clear the formatting
# again (ClearFormats) so the cell keeps plain text content without leaving
# a stray explicit style behind - this mirrors the original cells, which
# carry no style index at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("A2") "Turbine-A"
$ws.Range("B2").Value = 293.67
$ws.Range("C2").Value = 156.4
$ws.Range("D2").Value = 303.3
$ws.Range("E2").Value = 2403.29
$ws.Range("F2").Value = 25.18
Set-TextValue $ws.Range("G2") "86.2%"

# Row 3
Set-TextValue $ws.Range("A3") "Boiler-2"
$ws.Range("B3").Value = 187.05
$ws.Range("C3").Value = 147.76
$ws.Range("D3").Value = 307.69
$ws.Range("E3").Value = 3835.14
$ws.Range("F3").Value = 42.2
Set-TextValue $ws.Range("G3") "94.6%"

# Row 4
Set-TextValue $ws.Range("A4") "Turbine-A"
$ws.Range("B4").Value = 459.68
$ws.Range("C4").Value = 50.75
$ws.Range("D4").Value = 426.35
$ws.Range("E4").Value = 1795.1
$ws.Range("F4").Value = 29.14
Set-TextValue $ws.Range("G4") "88.7%"

# Row 5
Set-TextValue $ws.Range("A5") "Turbine-A"
$ws.Range("B5").Value = 447.51
$ws.Range("C5").Value = 84.32
$ws.Range("D5").Value = 449.87
$ws.Range("E5").Value = 2114.78
$ws.Range("F5").Value = 36.42
Set-TextValue $ws.Range("G5") "91.7%"

# Row 6
Set-TextValue $ws.Range("A6") "Cooling-Tower"
$ws.Range("B6").Value = 336.46
$ws.Range("C6").Value = 149.1
$ws.Range("D6").Value = 347.61
$ws.Range("E6").Value = 1985.8
$ws.Range("F6").Value = 26.09
Set-TextValue $ws.Range("G6") "89.5%"
